$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.481.13"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.690.33"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3886"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4027"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08755"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.525"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.988"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001350"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "1.686.05"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "98.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07097"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.265"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "24.483.07"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.970"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.352"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.769"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.223"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "1.871.51"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08844"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.446"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.037"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2827"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.961"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02925"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09129"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7957"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.454"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.619"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.207"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.352"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
